$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter Salvar em pdf Salvar em docx".
# The content footer that follows LOB1039's "Requisitos" entry consists of:
#   - an empty paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "(c) 2020 . Contact: ... Creative Commons Attribution"
# All three paragraphs (including the blank one right before "Ver no Jupiter")
# must be removed, while the paragraph mark of the preceding
# "LOB1039: ..." paragraph and the paragraph mark of the trailing blank
# paragraph (right before the page-break paragraph) must be preserved.

$verNoIndex = -1
$copyrightIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $verNoIndex = $idx
    }
    if ($p.Range.Text -like "*Powered by Jekyll*") {
        $copyrightIndex = $idx
    }
}

if ($verNoIndex -gt 0 -and $copyrightIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($verNoIndex - 1)
    $endPara = $d.Paragraphs.Item($copyrightIndex)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
